# TC05_Trials_Filter_Race-Unknown.xlsx : add a TabName/CasesTab column (col A)
# ahead of the existing query/StatQuery columns, and refresh the Cypher query
# text in the "query" (CasesTab) and "StatQuery" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = 'MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "UNKNOWN"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '''') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '''') AS `Trial Code`,
    COALESCE(a.arm_id, '''') AS `Arm`,
    COALESCE(a.arm_drug, '''') AS `Arm Treatment`,
    COALESCE(c.disease, '''') AS `Diagnosis`,
    COALESCE(c.gender, '''') AS `Gender`,
    COALESCE(c.race, '''') AS `Race`,
    COALESCE(c.ethnicity, '''') AS `Ethnicity`'

$statQuery = 'MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "UNKNOWN"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials'

# Insert a new column before column A, shifting everything right by one.
$ws.Range("A1").EntireColumn.Insert()

# Row 1 (headers)
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# Row 2 (values)
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = "TC05_Trials_Filter_Race-Unknown_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC05_Trials_Filter_Race-Unknown_WebData.xlsx"

# Wrap text on B2:C2 (query text cells) like the original query/StatQuery cells.
$ws.Range("B2:C2").WrapText = $true

# Column width for the new first column only - the other columns (B:E) keep
# the exact widths they already had (shifted right by the insert above).
$ws.Range("A1").EntireColumn.ColumnWidth = 8.81640625

# Row height for the data row (wrapped, taller text now)
$ws.Rows.Item(2).RowHeight = 174

# View / selection state
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("C5").Select()
